$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.233.51"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.840.25"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.39"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4645"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3876"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07873"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9652"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.14"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "1.863.79"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.706"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.899"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06874"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.68"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009977"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.76"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "28.241.95"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.313"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.113"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "2.052.52"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.70"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.20"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.748"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.977"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.24"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9382"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09270"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.294"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.328"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.332"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05835"
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02128"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.142"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.774"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5605"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.940"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1764"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07319"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.63"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5282"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.142"
$ws.Range("E46").Value = "  -8.61%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.141"
$ws.Range("E47").Value = "  -12.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.839"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.11"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.025"
$ws.Range("E51").Value = "  +0.32%  "
